$d = $word.ActiveDocument

$replacements = @(
    @("2025-03-28 Friday", "2025-03-29 Saturday"),
    @("325÷6=54, 1", "311÷2=155, 1"),
    @("740÷8=92, 4", "355÷4=88, 3"),
    @("498÷6=83, 0", "393÷6=65, 3"),
    @("889÷2=444, 1", "408÷4=102, 0"),
    @("555÷2=277, 1", "425÷2=212, 1"),
    @("297÷8=37, 1", "993÷5=198, 3"),
    @("728÷6=121, 2", "489÷8=61, 1"),
    @("672÷5=134, 2", "499÷9=55, 4"),
    @("278÷5=55, 3", "683÷3=227, 2"),
    @("867÷6=144, 3", "475÷4=118, 3"),
    @("432÷8=54, 0", "885÷2=442, 1"),
    @("800÷2=400, 0", "446÷4=111, 2"),
    @("482÷6=80, 2", "822÷2=411, 0"),
    @("557÷4=139, 1", "727÷9=80, 7"),
    @("262÷7=37, 3", "526÷8=65, 6"),
    @("622÷7=88, 6", "290÷2=145, 0"),
    @("457÷4=114, 1", "533÷7=76, 1"),
    @("849÷5=169, 4", "411÷8=51, 3"),
    @("319÷8=39, 7", "546÷9=60, 6"),
    @("958÷7=136, 6", "975÷8=121, 7"),
    @("526÷5=105, 1", "621÷6=103, 3"),
    @("252÷7=36, 0", "316÷2=158, 0"),
    @("683÷9=75, 8", "579÷7=82, 5"),
    @("598÷6=99, 4", "274÷4=68, 2"),
    @("421÷8=52, 5", "561÷8=70, 1")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
